$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 157, shifting existing rows 157:283 down to 158:284.
$ws.Rows.Item(157).Insert()

# Populate the newly inserted row 157 with the new data point.
$ws.Cells.Item(157, 1).Value = 8
$ws.Cells.Item(157, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(157, 3).Value = "Coquimbo"
$ws.Cells.Item(157, 4).Value = 44741
$ws.Cells.Item(157, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(157, 5).Value = 4
$ws.Cells.Item(157, 6).Value = 100112012
$ws.Cells.Item(157, 7).Value = "Espinaca"
$ws.Cells.Item(157, 8).Value = "Sin especificar"
$ws.Cells.Item(157, 9).Value = "Primera"
$ws.Cells.Item(157, 10).Value = 2700
$ws.Cells.Item(157, 11).Value = 500
$ws.Cells.Item(157, 12).Value = 600
$ws.Cells.Item(157, 13).Value = 550
$ws.Cells.Item(157, 14).Value = "$/atado 300 a 500 gramos"
$ws.Cells.Item(157, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(157, 16).Value = 1100
$ws.Cells.Item(157, 17).Value = 0.5
$ws.Cells.Item(157, 18).Value = "Hortaliza"
